$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @("D2", "96.830.62"),
    @("E2", "  +0.42%  "),
    @("D3", "3.650.50"),
    @("E3", "  +1.77%  "),
    @("E4", "  -0.01%  "),
    @("D5", "242.20"),
    @("E5", "  +0.15%  "),
    @("D6", "1.89"),
    @("E6", "  +21.46%  "),
    @("D7", "654.78"),
    @("E7", "  +0.03%  "),
    @("D8", "0.422"),
    @("E8", "  +3.98%  "),
    @("D9", "1.08"),
    @("E9", "  +4.25%  "),
    @("D11", "3.647.74"),
    @("E11", "  +1.82%  "),
    @("D12", "44.31"),
    @("E12", "  +2.64%  "),
    @("E13", "  +1.32%  "),
    @("D14", "6.50"),
    @("E14", "  +1.82%  "),
    @("D15", "4.330.69"),
    @("E15", "  +1.81%  "),
    @("D16", "96.690.06"),
    @("E16", "  +0.39%  "),
    @("D17", "0.0000258"),
    @("E17", "  -0.45%  "),
    @("D18", "3.647.74"),
    @("E18", "  +2.48%  "),
    @("B19", "Uniswap"),
    @("C19", "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"),
    @("D19", "12.87"),
    @("E19", "  +2.67%  "),
    @("B20", "Polkadot"),
    @("C20", "https://coinranking.com/coin/25W7FG7om+polkadot-dot"),
    @("D20", "7.76"),
    @("E20", "  +0.07%  "),
    @("D21", "18.32"),
    @("E21", "  +3.29%  "),
    @("D22", "0.537"),
    @("E22", "  +8.81%  "),
    @("D23", "511.85"),
    @("E23", "  +0.11%  "),
    @("D24", "3.43"),
    @("E24", "  -0.21%  "),
    @("D25", "0.0000205"),
    @("E25", "  +1.29%  "),
    @("E26", "  +0.92%  "),
    @("D27", "101.29"),
    @("E27", "  +5.10%  "),
    @("D28", "13.06"),
    @("E28", "  +2.35%  "),
    @("E29", "  +15.77%  "),
    @("D30", "3.02"),
    @("E30", "  +1.37%  "),
    @("D31", "11.90"),
    @("E31", "  +4.06%  "),
    @("D32", "0.999"),
    @("E32", "  -0.07%  "),
    @("E33", "  +1.17%  "),
    @("D34", "32.99"),
    @("E34", "  +4.20%  "),
    @("D35", "0.999"),
    @("E35", "  -0.07%  "),
    @("D36", "1.74"),
    @("E36", "  +7.26%  "),
    @("D37", "0.586"),
    @("E37", "  +3.78%  "),
    @("B38", "Bittensor"),
    @("C38", "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"),
    @("D38", "613.94"),
    @("E38", "  +0.20%  "),
    @("B39", "RenderToken"),
    @("C39", "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"),
    @("D39", "8.74"),
    @("E39", "  +1.19%  "),
    @("D40", "41.37"),
    @("E40", "  +20.57%  "),
    @("D41", "0.160"),
    @("E41", "  +5.90%  "),
    @("D42", "0.947"),
    @("E42", "  +4.35%  "),
    @("D43", "1.92"),
    @("E43", "  +5.67%  "),
    @("D45", "6.17"),
    @("E45", "  +8.29%  "),
    @("D46", "0.0446"),
    @("E46", "  +7.06%  "),
    @("D47", "0.430"),
    @("E47", "  +28.15%  "),
    @("D48", "2.30"),
    @("E48", "  +1.20%  "),
    @("E49", "  +0.33%  "),
    @("D50", "8.64"),
    @("E50", "  +5.97%  "),
    @("D51", "54.39"),
    @("E51", "  +2.19%  "),
)

foreach ($pair in $updates) {
    $cellRef = $pair[0]
    $newVal = $pair[1]
    $rng = $ws.Range($cellRef)
    if ($cellRef.StartsWith("D")) {
        # Force text storage so numeric-looking strings (e.g. "242.20") are not
        # converted into floating point numbers, matching the original inlineStr type.
        $rng.NumberFormat = "@"
        $rng.Value = $newVal
        $rng.Style = "Normal"
    } else {
        $rng.Value = $newVal
    }
}
